# Update the "想去人数" (want-to-go count) values in column F on both the
# "展览" and "全部类型" worksheets, which carry duplicate listings of the
# same convention events.

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    3  = 1072
    4  = 126
    7  = 54
    8  = 11117
    9  = 4267
    13 = 2494
    15 = 87
    17 = 153
    18 = 480
    19 = 11209
    20 = 11052
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
